# Applies the "scenario 3b_1" coverage update described by the commit
# "add scenario files for runs march 8th 2023":
#   - Row 2 (New Product A, All Treatment Campaign MDA, 5-15yo): flat coverage
#     bumped from 0.6 to 0.736 across every year column (H:AD).
#   - Rows 6,7,9: a handful of near-term year cells (previously 0, still
#     ramping up) get non-zero coverage values.
#   - Rows 8,10: the T:V cells lose their quote-prefixed "text number" style
#     and pick up real coverage values (row 8) / stay 0 (row 10 U:V), while
#     T10 also gets a value.
#   - Row 11: T11 moves from 0 to 0.45.
#   - A brand new row 12 ("Vector Control") is appended with two tiny
#     (1e-8) values in J12/K12, using the same style as the other
#     quote-prefixed numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Platform Coverage")

# --- Row 2: H2:AD2 all become 0.736 -----------------------------------
$ws.Range("H2:AD2").Value = 0.736

# --- Row 6: T6:W6 ramp up -----------------------------------------------
$ws.Range("T6").Value = 0.79
$ws.Range("U6").Value = 0.8
$ws.Range("V6").Value = 0.81
$ws.Range("W6").Value = 0.82

# --- Row 7: U7:W7 ramp up -------------------------------------------------
$ws.Range("U7").Value = 0.68
$ws.Range("V7").Value = 0.68849999999999989
$ws.Range("W7").Value = 0.69700000000000006

# --- Row 8: T8:V8 lose the quote-prefixed style; T8 gets a real value ---
$ws.Range("T8").Value = 0.45
$ws.Range("U8").Value = 0
$ws.Range("V8").Value = 0
$ws.Range("T8:V8").Style = "Normal"

# --- Row 9: U9:W9 ramp up (same pattern as row 7) ------------------------
$ws.Range("U9").Value = 0.68
$ws.Range("V9").Value = 0.68849999999999989
$ws.Range("W9").Value = 0.69700000000000006

# --- Row 10: T10:V10 lose the quote-prefixed style; T10 gets a value ----
$ws.Range("T10").Value = 0.6715000000000001
$ws.Range("U10").Value = 0
$ws.Range("V10").Value = 0
$ws.Range("T10:V10").Style = "Normal"

# --- Row 11: T11 ----------------------------------------------------------
$ws.Range("T11").Value = 0.45

# --- Row 12 (new): Vector Control ----------------------------------------
$ws.Range("B12").Value = "Vector Control"
$ws.Range("J12").Value = 0.00000001
$ws.Range("K12").Value = 0.00000001
$ws.Range("J12:K12").Font.Color = 0

# --- Selection: matches the saved cursor position in the target file -----
[void]$ws.Range("I2:AD2").Select()
